$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Added test case "CMC Autopay 2.7" - three new sheets copied from an
# existing PayNow* template sheet (same 14-column schema used by every
# sheet in this workbook), appended after CCDeferredCorp_27, with their
# Result/Date/EmulatorData/NameID/ACHID data filled in.
# ---------------------------------------------------------------------------

$template = $wb.Worksheets.Item("PayNowNoCFPS_27")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- CMCAutopayPC_27 --------------------------------------------------
$template.Copy($null, $lastSheet)
$wsPC = $wb.Worksheets.Item($lastSheet.Index + 1)
$wsPC.Name = "CMCAutopayPC_27"
$wsPC.Range("B2").Value = "Wed Jan 29 17:54:11 IST 2025"
$wsPC.Range("E2").Value = "19"
$wsPC.Range("K2").Value = "4"
$wsPC.Range("H2").ClearContents() | Out-Null
$wsPC.Range("D2:H2").Select() | Out-Null

# --- CMCAutopayCorp_27 -------------------------------------------------
$template.Copy($null, $wsPC)
$wsCorp = $wb.Worksheets.Item($wsPC.Index + 1)
$wsCorp.Name = "CMCAutopayCorp_27"
$wsCorp.Range("B2").Value = "Wed Feb 05 17:20:04 IST 2025"
$wsCorp.Range("E2").Value = "19"
$wsCorp.Range("J2").Value = "3"
$wsCorp.Range("K2").Value = "3"
$wsCorp.Range("H2").ClearContents() | Out-Null
$wsCorp.Range("D2:H2").Select() | Out-Null

# --- CMCAutopayPS_27 ----------------------------------------------------
$template.Copy($null, $wsCorp)
$wsPS = $wb.Worksheets.Item($wsCorp.Index + 1)
$wsPS.Name = "CMCAutopayPS_27"
$wsPS.Range("B2").Value = "Tue Feb 04 19:16:37 IST 2025"
$wsPS.Range("E2").Value = "19"
$wsPS.Range("H2").ClearContents() | Out-Null
$wsPS.Range("M8").Select() | Out-Null

# An unrelated navigation tweak on PayNowNoCFPS_27 captured by the same save.
$template.Range("A1:N1").Select() | Out-Null

# Previously-last sheet no longer holds the "selected tab" -- the new last
# sheet (CMCAutopayPS_27) becomes the active / tabSelected sheet instead.
$wsPS.Activate() | Out-Null
